$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove the "correo" / "fono" (email / phone) columns' data ---
# Header row: keep H1 as an empty (still-styled) cell, remove I1 entirely.
$ws.Range("H1").ClearContents()
$ws.Range("I1").ClearContents()
$ws.Range("I1").ClearFormats()

# Body rows 2-11: clear correo/fono values but keep their formatting.
$ws.Range("H2:H11").ClearContents()
$ws.Range("I2:I11").ClearContents()

# Last data row (12): clear correo and remove fono entirely (row ends at H).
$ws.Range("H12").ClearContents()
$ws.Range("I12").ClearContents()
$ws.Range("I12").ClearFormats()

# --- Remove the hyperlink that was attached to H3 ---
$ws.Range("H3").Hyperlinks.Delete()

# --- Append two blank rows (13 and 14) below the data, matching body style ---
$ws.Range("A13:H14").Value = ""
